# Generate Report for Archive
#
# The localization status for the handed-off items moved from
# "Ready for handoff" to "In Translation" everywhere it appears
# (Overview!E2:F4, zh-cn!C2:C4, de-de!C2:C4), and the "Status" columns
# were re-fit to the new (shorter) text.

$wb = $excel.ActiveWorkbook

# Update the status text on every sheet wherever it occurs.
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# Re-fit the status columns now that the text is shorter.
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5
